# Applies the content edits found on slide 11 of the STAMP SRPM deck:
#   1. Title text rewording + font size 30pt -> 28pt
#   2. "Sender:" -> "Session-Sender:"
#   3. "Reflector:" -> "Session-Reflector:"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# --- 1. Title: reword + shrink font from 30pt to 28pt ---------------------
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "STAMP - Stand-alone Direct-mode LM Test Packet Format"
$titleRange.Font.Size = 28

# --- 2/3. Content placeholder: "Sender:" / "Reflector:" relabel -----------
$content = $s.Shapes.Item(4)
$contentRange = $content.TextFrame.TextRange
$fullText = $contentRange.Text

$senderIdx = $fullText.IndexOf("Sender:")
$contentRange.Characters($senderIdx + 1, 7).Text = "Session-Sender:"

# Text shifted after the previous edit, so re-read before locating "Reflector:"
$fullText = $contentRange.Text
$reflectorIdx = $fullText.IndexOf("Reflector:")
$contentRange.Characters($reflectorIdx + 1, 10).Text = "Session-Reflector:"
